$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns remain text so numeric-looking strings are not reinterpreted as numbers
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.981.07"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.640.46"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.66%  "
$ws.Range("D5").Value = "215.06"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "0.5097"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").Value = "0.2585"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "0.06357"
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("D10").Value = "19.83"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("D11").Value = "0.07765"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "4.278"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("D13").Value = "1.633.56"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").Value = "0.5470"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "0.0₅7749"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").Value = "64.34"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "26.002.95"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "196.11"
$ws.Range("E19").Value = "  -1.47%  "
$ws.Range("D20").Value = "4.433"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").Value = "9.927"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").Value = "6.094"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "143.51"
$ws.Range("E25").Value = "  +2.16%  "
$ws.Range("D26").Value = "0.1236"
$ws.Range("E26").Value = "  +7.26%  "
$ws.Range("D27").Value = "6.864"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("D28").Value = "15.62"
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").Value = "0.04862"
$ws.Range("E30").Value = "  -3.63%  "
$ws.Range("D31").Value = "3.273"
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "3.228"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "1.543"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").Value = "2.379"
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").Value = "0.9134"
$ws.Range("E35").Value = "  +1.83%  "
$ws.Range("D36").Value = "2.570"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").Value = "0.5549"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Value = "1.092.34"
$ws.Range("E38").Value = "  -4.20%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "1.002"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").Value = "2.528"
$ws.Range("E41").Value = "  -1.53%  "
$ws.Range("D42").Value = "5.594"
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("D43").Value = "0.8059"
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("D44").Value = "99.18"
$ws.Range("E45").Value = "  -3.88%  "
$ws.Range("D46").Value = "1.779.54"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").Value = "0.4536"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "1.009"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("D49").Value = "55.26"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").Value = "0.05206"
$ws.Range("E50").Value = "  +2.16%  "
$ws.Range("D51").Value = "7.488"
$ws.Range("E51").Value = "  +0.97%  "
